# Apply the edit described by the diff:
#  - Insert a new column F (former column E shifts right, keeping all of its
#    data/formatting) so the sheet grows from A:E to A:F
#  - Header E1 stays "Preprocesamiento 2"; new header F1 = "Preprocesamiento 3"
#  - New column E (rows 2-24) gets "Analizar outliers" only where
#    Preprocesamiento 2 (col D) = "Imputar faltantes (como?)" AND the old
#    column E value (now sitting in col F) was "Normalizar"; otherwise blank
#  - Column E width matches column D's former width (27.28515625)
#  - Update the sheet view (topLeftCell A5 / selection E26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; everything that was in E (data + formats) moves to F
$ws.Range("E1").EntireColumn.Insert()

# Restore the header text (Insert pushed "Preprocesamiento 2" into F1) and
# give the new column its own header.
$ws.Range("E1").Value = "Preprocesamiento 2"
$ws.Range("F1").Value = "Preprocesamiento 3"

# Fill the new column E: "Analizar outliers" where Preprocesamiento 2 is
# "Imputar faltantes (como?)" and Preprocesamiento 3 (old E, now F) is
# "Normalizar"; blank everywhere else.
for ($r = 2; $r -le 24; $r++) {
    $prep2 = $ws.Cells.Item($r, 4).Value   # column D
    $prep3 = $ws.Cells.Item($r, 6).Value   # column F (former E)

    if ($prep2 -eq "Imputar faltantes (como?)" -and $prep3 -eq "Normalizar") {
        $ws.Cells.Item($r, 5).Value = "Analizar outliers"
    } else {
        $ws.Cells.Item($r, 5).Value = $null
    }
}

# Column widths: E keeps the width D used to have; F keeps the width old E had.
$ws.Columns.Item(5).ColumnWidth = 27.28515625
$ws.Columns.Item(6).ColumnWidth = 19.85546875

# Update sheet view: scroll so row 5 is the top row, and select E26.
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("E26").Select()
